$wb = $excel.ActiveWorkbook

# --- Clear header styling (A1:N1) on all sheets: bold/border/center-align -> default ---
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").ClearFormats()
}

# --- Data value updates (market price refresh) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115.7
$ws.Range("I9").Value = 113.6
$ws.Range("J9").Value = 117.8
$ws.Range("K9").Value = 113.6
$ws.Range("L9").Value = 117.8
$ws.Range("M9").Value = 55.40000000000001
$ws.Range("N9").Value = -455.8
$ws.Range("H38").Value = 3458
$ws.Range("H40").Value = 3849.92
$ws.Range("J40").Value = 4341.0586
$ws.Range("L40").Value = 4341.0586
$ws.Range("N40").Value = -4691.0586
$ws.Range("H51").Value = 45138
$ws.Range("J51").Value = 45138
$ws.Range("L51").Value = 45138
$ws.Range("N51").Value = -46106
$ws.Range("H58").Value = 17668.834
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 17668.834
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 53006.50199999999
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -53306.50199999999
$ws.Range("H62").Value = 5146.125
$ws.Range("I62").Value = 3848.9092
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 3848.9092
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -3224.9092
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 5146.125
$ws.Range("I65").Value = 3848.9092
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 19244.546
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -16124.546
$ws.Range("N65").Value = -46240
$ws.Range("H70").Value = 3231.3333
$ws.Range("I70").Value = 1195
$ws.Range("J70").Value = 4249.5
$ws.Range("K70").Value = 3585
$ws.Range("L70").Value = 12748.5
$ws.Range("M70").Value = -3315
$ws.Range("N70").Value = -13288.5
$ws.Range("H73").Value = 3231.3333
$ws.Range("I73").Value = 1195
$ws.Range("J73").Value = 4249.5
$ws.Range("K73").Value = 3585
$ws.Range("L73").Value = 12748.5
$ws.Range("M73").Value = -2649
$ws.Range("N73").Value = -14620.5
$ws.Range("H86").Value = 7266.923
$ws.Range("I86").Value = 12499.667
$ws.Range("K86").Value = 12499.667
$ws.Range("M86").Value = -11376.667
$ws.Range("H89").Value = 7266.923
$ws.Range("I89").Value = 12499.667
$ws.Range("K89").Value = 62498.335
$ws.Range("M89").Value = -56882.335
$ws.Range("H132").Value = 1888.4906
$ws.Range("I132").Value = 1523.6666
$ws.Range("K132").Value = 4570.9998
$ws.Range("M132").Value = -2040.9998
$ws.Range("H138").Value = 3256.0173
$ws.Range("I138").Value = 2110.4211
$ws.Range("K138").Value = 6331.263300000001
$ws.Range("M138").Value = -1191.263300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7584077.5
$ws.Range("I32").Value = 8200504.5
$ws.Range("K32").Value = 8200504.5
$ws.Range("M32").Value = -8200217.5
$ws.Range("H74").Value = 6255564.5
$ws.Range("I74").Value = 8065758.5
$ws.Range("J74").Value = 20451.555
$ws.Range("K74").Value = 8065758.5
$ws.Range("L74").Value = 20451.555
$ws.Range("M74").Value = -8064884.5
$ws.Range("N74").Value = -22199.555
$ws.Range("H77").Value = 6255564.5
$ws.Range("I77").Value = 8065758.5
$ws.Range("J77").Value = 20451.555
$ws.Range("K77").Value = 40328792.5
$ws.Range("L77").Value = 102257.775
$ws.Range("M77").Value = -40324424.5
$ws.Range("N77").Value = -110993.775
$ws.Range("H95").Value = 46735.668
$ws.Range("J95").Value = 46735.668
$ws.Range("L95").Value = 46735.668
$ws.Range("N95").Value = -52227.668
$ws.Range("H132").Value = 5369.8
$ws.Range("I132").Value = 2215.8057
$ws.Range("K132").Value = 6647.4171
$ws.Range("M132").Value = -4117.4171

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 129995
$ws.Range("J82").Value = 129995
$ws.Range("L82").Value = 129995
$ws.Range("N82").Value = -130761
$ws.Range("H85").Value = 129995
$ws.Range("J85").Value = 129995
$ws.Range("L85").Value = 129995
$ws.Range("N85").Value = -132647
$ws.Range("H134").Value = 25110.488
$ws.Range("I134").Value = 2710.639
$ws.Range("K134").Value = 8131.917
$ws.Range("M134").Value = -5596.917

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1957
$ws.Range("I58").Value = 1696.25
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1696.25
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1493.25
$ws.Range("N58").Value = -3406
$ws.Range("H74").Value = 153571.42
$ws.Range("J74").Value = 153571.42
$ws.Range("L74").Value = 153571.42
$ws.Range("N74").Value = -155319.42
$ws.Range("H77").Value = 153571.42
$ws.Range("J77").Value = 153571.42
$ws.Range("L77").Value = 460714.26
$ws.Range("N77").Value = -469450.26
$ws.Range("H107").Value = 754.1111
$ws.Range("I107").Value = 569.5714
$ws.Range("K107").Value = 569.5714
$ws.Range("M107").Value = 1350.4286
$ws.Range("H136").Value = 1957
$ws.Range("I136").Value = 1696.25
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5088.75
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2538.75
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 427.64285
$ws.Range("I14").Value = 427.64285
$ws.Range("K14").Value = 1282.92855
$ws.Range("M14").Value = -1109.92855
$ws.Range("H49").Value = 3250.25
$ws.Range("I49").Value = 999
$ws.Range("K49").Value = 2997
$ws.Range("M49").Value = -2841
$ws.Range("H56").Value = 7683.381
$ws.Range("I56").Value = 7683.381
$ws.Range("K56").Value = 7683.381
$ws.Range("M56").Value = -7153.381
$ws.Range("H68").Value = 35185.7
$ws.Range("J68").Value = 40291.23
$ws.Range("L68").Value = 120873.69
$ws.Range("N68").Value = -122495.69
$ws.Range("H70").Value = 6251.5
$ws.Range("I70").Value = 4877.25
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 14631.75
$ws.Range("L70").Value = 27000
$ws.Range("M70").Value = -14316.75
$ws.Range("N70").Value = -27630
$ws.Range("H71").Value = 35185.7
$ws.Range("J71").Value = 40291.23
$ws.Range("L71").Value = 362621.07
$ws.Range("N71").Value = -370733.07
$ws.Range("H73").Value = 6251.5
$ws.Range("I73").Value = 4877.25
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 14631.75
$ws.Range("L73").Value = 27000
$ws.Range("M73").Value = -13539.75
$ws.Range("N73").Value = -29184
$ws.Range("H87").Value = 6486.857
$ws.Range("I87").Value = 3401.3333
$ws.Range("K87").Value = 10203.9999
$ws.Range("M87").Value = -8955.999899999999
$ws.Range("H90").Value = 6486.857
$ws.Range("I90").Value = 3401.3333
$ws.Range("K90").Value = 30611.9997
$ws.Range("M90").Value = -24371.9997
$ws.Range("H107").Value = 665.4286
$ws.Range("J107").Value = 1110.5
$ws.Range("L107").Value = 3331.5
$ws.Range("N107").Value = -7171.5
$ws.Range("H110").Value = 14999
$ws.Range("I110").Value = 14999
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 44997
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -40907
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 539.37933
$ws.Range("I122").Value = 456.07693
$ws.Range("J122").Value = 607.0625
$ws.Range("K122").Value = 4104.69237
$ws.Range("L122").Value = 5463.5625
$ws.Range("M122").Value = -1654.69237
$ws.Range("N122").Value = -10363.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 17254.25
$ws.Range("I43").Value = 9672.333
$ws.Range("J43").Value = 40000
$ws.Range("K43").Value = 9672.333
$ws.Range("L43").Value = 40000
$ws.Range("M43").Value = -9521.333
$ws.Range("N43").Value = -40302
$ws.Range("H102").Value = 2594.0588
$ws.Range("I102").Value = 1927.9584
$ws.Range("K102").Value = 1927.9584
$ws.Range("M102").Value = -305.9584
$ws.Range("H126").Value = 4507
$ws.Range("I126").Value = 4504
$ws.Range("K126").Value = 13512
$ws.Range("M126").Value = -11042
$ws.Range("H132").Value = 52635024
$ws.Range("I132").Value = 62503452
$ws.Range("J132").Value = 3415.3333
$ws.Range("K132").Value = 187510356
$ws.Range("L132").Value = 10245.9999
$ws.Range("M132").Value = -187507826
$ws.Range("N132").Value = -15305.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 104170.8
$ws.Range("I7").Value = 3538.8
$ws.Range("K7").Value = 3538.8
$ws.Range("M7").Value = -3426.8
$ws.Range("H22").Value = 8699.4
$ws.Range("I22").Value = 8374.25
$ws.Range("K22").Value = 8374.25
$ws.Range("M22").Value = -8079.25
$ws.Range("H27").Value = 8699.4
$ws.Range("I27").Value = 8374.25
$ws.Range("K27").Value = 8374.25
$ws.Range("M27").Value = -8267.25
$ws.Range("H40").Value = 4180.7334
$ws.Range("I40").Value = 3244
$ws.Range("K40").Value = 3244
$ws.Range("M40").Value = -3108
$ws.Range("H55").Value = 55556056
$ws.Range("I55").Value = 90909704
$ws.Range("J55").Value = 318.85715
$ws.Range("K55").Value = 90909704
$ws.Range("L55").Value = 318.85715
$ws.Range("M55").Value = -90909531
$ws.Range("N55").Value = -664.85715
$ws.Range("H82").Value = 1892.5714
$ws.Range("I82").Value = 1892.5714
$ws.Range("K82").Value = 1892.5714
$ws.Range("M82").Value = -1531.5714
$ws.Range("H85").Value = 1892.5714
$ws.Range("I85").Value = 1892.5714
$ws.Range("K85").Value = 1892.5714
$ws.Range("M85").Value = -644.5714
$ws.Range("H101").Value = 9997.2
$ws.Range("J101").Value = 9997.2
$ws.Range("L101").Value = 9997.2
$ws.Range("N101").Value = -16487.2
$ws.Range("H126").Value = 104170.8
$ws.Range("I126").Value = 3538.8
$ws.Range("K126").Value = 10616.4
$ws.Range("M126").Value = -8146.400000000001
$ws.Range("H136").Value = 127086
$ws.Range("I136").Value = 88004
$ws.Range("J136").Value = 205250
$ws.Range("K136").Value = 264012
$ws.Range("L136").Value = 615750
$ws.Range("M136").Value = -261462
$ws.Range("N136").Value = -620850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3993.8298
$ws.Range("I126").Value = 3291.4243
$ws.Range("J126").Value = 5649.5
$ws.Range("K126").Value = 9874.2729
$ws.Range("L126").Value = 16948.5
$ws.Range("M126").Value = -7404.2729
$ws.Range("N126").Value = -21888.5
$ws.Range("H132").Value = 281183.1
$ws.Range("I132").Value = 1882.3793
$ws.Range("K132").Value = 5647.1379
$ws.Range("M132").Value = -3117.1379

